# feat: add 2022-Q3 data
#
# Before: sheets "总计", "2022-Q1"
# After:  sheets "总计", "2022-Q3", "2022-Q1"  (new "2022-Q3" sheet inserted
#         between the two, a new summary row added on "总计", and the old
#         "2022-Q1" detail sheet duplicated so a separate "2022-Q1" sheet
#         keeps its original data alongside the new "2022-Q3" data).

$wb = $excel.ActiveWorkbook
$total = $wb.Worksheets.Item("总计")

# ---------------------------------------------------------------------------
# 1. Duplicate the existing "2022-Q1" worksheet. The original keeps its
#    position/sheetId and becomes "2022-Q3" (new data); the duplicate is
#    placed right after it and keeps the "2022-Q1" name/data untouched.
# ---------------------------------------------------------------------------
$q1 = $wb.Worksheets.Item("2022-Q1")
$q1.Copy($null, $q1)

$q3 = $q1
$q3.Name = "2022-Q3"

$q1New = $wb.Worksheets.Item("2022-Q1 (2)")
$q1New.Name = "2022-Q1"

# ---------------------------------------------------------------------------
# 2. "2022-Q3" only has one fund row (vs. two on the old "2022-Q1" sheet) -
#    drop row 3 and overwrite row 2 with the new fund's data. The header row
#    and the A-column counter on "2022-Q3" use the same cell style as the
#    "总计" sheet's header/counter (style index 2), not the style the old
#    "2022-Q1" sheet used (style index 1) - match that explicitly.
# ---------------------------------------------------------------------------
$q3.Rows.Item(3).Delete()

$total.Range("B1").Copy()
$q3.Range("B1:H1").PasteSpecial(-4122)
$total.Range("A2").Copy()
$q3.Range("A2").PasteSpecial(-4122)

$q3.Range("B2:G2").NumberFormat = "@"
$q3.Range("B2").Value = "004448"
$q3.Range("C2").Value = "博时汇智回报灵活配置混合"
$q3.Range("D2").Value = "1.77"
$q3.Range("E2").Value = "67.69"
$q3.Range("F2").Value = "3.16"
$q3.Range("G2").Value = "0.0559"
$q3.Range("H2").Value = 9

# ---------------------------------------------------------------------------
# 3. Update the "总计" summary sheet: push the existing "2022-Q1" row down to
#    row 3 and write the new "2022-Q3" totals into row 2.
# ---------------------------------------------------------------------------
$oldB = $total.Range("B2").Value2
$oldC = $total.Range("C2").Value2
$oldD = $total.Range("D2").Value2

$total.Range("A2").Copy()
$total.Range("A3").PasteSpecial(-4122)

$total.Range("A3").Value = 1
$total.Range("B3").Value = $oldB
$total.Range("C3").Value = $oldC
$total.Range("D3").Value = $oldD

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q3"
$total.Range("C2").Value = 1
$total.Range("D2").Value = 0.06
